# Implemented doorbell chime selection and began video capability
#
# 1. Remove the stray orphaned URL that lived in Main Icons!F5 (and thus
#    the now-unused shared string).
# 2. Add a new "Videos" sheet right after "Main Icons" with three rows
#    of sample video links (Grandma, Cool Guy, Couple), and make it the
#    active/selected tab.
# 3. Append two new rows to the Main Icons sheet: Camera, Chat.

$wb = $excel.ActiveWorkbook

# --- Main Icons sheet -------------------------------------------------
$mainIcons = $wb.Worksheets.Item("Main Icons")

# Drop the orphaned home-security.png reference that lived outside the
# A:B table in F5.
$mainIcons.Range("F5").Clear()

# --- New Videos sheet --------------------------------------------------
$videos = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $mainIcons)
$videos.Name = "Videos"

$videos.Range("A1").Value = "Grandma"
$videos.Range("B1").Value = "http://www.shutterstock.com/video/clip-7209262-stock-footage-silly-grandma-making-funny-faces-at-the-camera.html"

# New doorbell-chime-selection related rows (Camera / Chat) appended
# after the existing "Language Flags" row.
$mainIcons.Range("A31").Value = "Camera"
$mainIcons.Range("B31").Value = "http://www.flaticon.com/free-icon/photo-camera_3901"
$mainIcons.Range("A32").Value = "Chat"
$mainIcons.Range("B32").Value = "http://www.flaticon.com/free-icon/chat-bubbles-with-ellipsis_61516"

$mainIcons.Range("B32").Select()

$videos.Range("A2").Value = "Cool Guy"
$videos.Range("B2").Value = "https://www.shutterstock.com/video/clip-11529716-stock-footage-handsome-man-in-a-light-blue-jacket-standing-outside-on-a-sunny-summer-day-with-a-green-background.html"
$videos.Range("A3").Value = "Couple"
$videos.Range("B3").Value = "http://www.shutterstock.com/video/clip-5574770-stock-footage-young-happy-couple-standing-on-house-porch.html"

$videos.Range("B3").Select()
$videos.Activate()
